$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.630712666666667
$ws.Range("H2").Value = 7.892138000000001
$ws.Range("I2").Value = 0.3947434022685045
$ws.Range("J2").Value = 0.3947434022685045
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1199136666666667
$ws.Range("N2").Value = 0.359741
$ws.Range("O2").Value = 0.03475238343556118
$ws.Range("P2").Value = 0.03475238343556118
$ws.Range("Q2").Value = 0.3154584018064445
$ws.Range("R2").Value = 2.839125616258
$ws.Range("S2").Value = 0.01371827407429304
$ws.Range("T2").Value = 0.01371827407429304
$ws.Range("G3").Value = 2.630712666666667
$ws.Range("H3").Value = 7.892138000000001
$ws.Range("I3").Value = 0.3947434022685045
$ws.Range("J3").Value = 0.3947434022685045
$ws.Range("N3").Value = 0.7418130000000001
$ws.Range("O3").Value = 0.07166202855244176
$ws.Range("P3").Value = 0.07166202855244176
$ws.Range("Q3").Value = 0.6504989517993335
$ws.Range("R3").Value = 5.854490566194001
$ws.Range("S3").Value = 0.02828811296425357
$ws.Range("T3").Value = 0.02828811296425357
$ws.Range("G4").Value = 2.630712666666667
$ws.Range("H4").Value = 7.892138000000001
$ws.Range("I4").Value = 0.3947434022685045
$ws.Range("J4").Value = 0.3947434022685045
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.355759
$ws.Range("N4").Value = 1.067277
$ws.Range("O4").Value = 0.1031031201224087
$ws.Range("P4").Value = 0.1031031201224087
$ws.Range("Q4").Value = 0.9358997075806668
$ws.Range("R4").Value = 8.423097368226001
$ws.Range("S4").Value = 0.04069927642161793
$ws.Range("T4").Value = 0.04069927642161793
$ws.Range("G5").Value = 2.630712666666667
$ws.Range("H5").Value = 7.892138000000001
$ws.Range("I5").Value = 0.3947434022685045
$ws.Range("J5").Value = 0.3947434022685045
$ws.Range("M5").Value = 2.727572666666667
$ws.Range("N5").Value = 8.182718000000001
$ws.Range("O5").Value = 0.7904824678895883
$ws.Range("P5").Value = 0.7904824678895883
$ws.Range("Q5").Value = 7.17545996345378
$ws.Range("R5").Value = 64.57913967108402
$ws.Range("S5").Value = 0.31203773880834
$ws.Range("T5").Value = 0.31203773880834
$ws.Range("I6").Value = 0.5991736942634763
$ws.Range("J6").Value = 0.5991736942634763
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1199136666666667
$ws.Range("N6").Value = 0.359741
$ws.Range("O6").Value = 0.03475238343556118
$ws.Range("P6").Value = 0.03475238343556118
$ws.Range("Q6").Value = 0.4788284615033333
$ws.Range("R6").Value = 4.30945615353
$ws.Range("S6").Value = 0.02082271396754604
$ws.Range("T6").Value = 0.02082271396754604
$ws.Range("I7").Value = 0.5991736942634763
$ws.Range("J7").Value = 0.5991736942634763
$ws.Range("N7").Value = 0.7418130000000001
$ws.Range("O7").Value = 0.07166202855244176
$ws.Range("P7").Value = 0.07166202855244176
$ws.Range("Q7").Value = 0.9873803028100001
$ws.Range("R7").Value = 8.886422725290002
$ws.Range("S7").Value = 0.04293800238618125
$ws.Range("T7").Value = 0.04293800238618125
$ws.Range("I8").Value = 0.5991736942634763
$ws.Range("J8").Value = 0.5991736942634763
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.355759
$ws.Range("N8").Value = 1.067277
$ws.Range("O8").Value = 0.1031031201224087
$ws.Range("P8").Value = 0.1031031201224087
$ws.Range("Q8").Value = 1.42058482049
$ws.Range("R8").Value = 12.78526338441
$ws.Range("S8").Value = 0.06177667737383459
$ws.Range("T8").Value = 0.06177667737383459
$ws.Range("I9").Value = 0.5991736942634763
$ws.Range("J9").Value = 0.5991736942634763
$ws.Range("M9").Value = 2.727572666666667
$ws.Range("N9").Value = 8.182718000000001
$ws.Range("O9").Value = 0.7904824678895883
$ws.Range("P9").Value = 0.7904824678895883
$ws.Range("Q9").Value = 10.89149769099333
$ws.Range("R9").Value = 98.02347921894003
$ws.Range("S9").Value = 0.4736363005359144
$ws.Range("T9").Value = 0.4736363005359144
$ws.Range("G10").Value = 0.04053866666666667
$ws.Range("H10").Value = 0.121616
$ws.Range("I10").Value = 0.00608290346801924
$ws.Range("J10").Value = 0.006082903468019241
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.1199136666666667
$ws.Range("N10").Value = 0.359741
$ws.Range("O10").Value = 0.03475238343556118
$ws.Range("P10").Value = 0.03475238343556118
$ws.Range("Q10").Value = 0.004861140161777778
$ws.Range("R10").Value = 0.043750261456
$ws.Range("S10").Value = 0.0002113953937221095
$ws.Range("T10").Value = 0.0002113953937221095
$ws.Range("G11").Value = 0.04053866666666667
$ws.Range("H11").Value = 0.121616
$ws.Range("I11").Value = 0.00608290346801924
$ws.Range("J11").Value = 0.006082903468019241
$ws.Range("N11").Value = 0.7418130000000001
$ws.Range("O11").Value = 0.07166202855244176
$ws.Range("P11").Value = 0.07166202855244176
$ws.Range("Q11").Value = 0.01002403664533333
$ws.Range("R11").Value = 0.09021632980800001
$ws.Range("S11").Value = 0.0004359132020069418
$ws.Range("T11").Value = 0.0004359132020069419
$ws.Range("G12").Value = 0.04053866666666667
$ws.Range("H12").Value = 0.121616
$ws.Range("I12").Value = 0.00608290346801924
$ws.Range("J12").Value = 0.006082903468019241
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.355759
$ws.Range("N12").Value = 1.067277
$ws.Range("O12").Value = 0.1031031201224087
$ws.Range("P12").Value = 0.1031031201224087
$ws.Range("Q12").Value = 0.01442199551466667
$ws.Range("R12").Value = 0.129797959632
$ws.Range("S12").Value = 0.0006271663269562043
$ws.Range("T12").Value = 0.0006271663269562044
$ws.Range("G13").Value = 0.04053866666666667
$ws.Range("H13").Value = 0.121616
$ws.Range("I13").Value = 0.00608290346801924
$ws.Range("J13").Value = 0.006082903468019241
$ws.Range("M13").Value = 2.727572666666667
$ws.Range("N13").Value = 8.182718000000001
$ws.Range("O13").Value = 0.7904824678895883
$ws.Range("P13").Value = 0.7904824678895883
$ws.Range("Q13").Value = 0.1105721591431111
$ws.Range("R13").Value = 0.9951494322880001
$ws.Range("S13").Value = 0.004808428545333985
$ws.Range("T13").Value = 0.004808428545333985
